$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings are preserved exactly
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.401.82"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "1.842.63"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "239.47"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").Value = "0.6255"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("D8").Value = "0.07436"
$ws.Range("E8").Value = "  -0.79%  "

$ws.Range("E9").Value = "  +2.34%  "

$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("D11").Value = "0.07724"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").Value = "1.842.95"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("D14").Value = "0.6756"
$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("E15").Value = "  -2.28%  "

$ws.Range("D16").Value = "81.92"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").Value = "6.233"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").Value = "29.399.01"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").Value = "233.35"
$ws.Range("E19").Value = "  +1.92%  "

$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D22").Value = "7.323"
$ws.Range("E22").Value = "  -2.19%  "

$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").Value = "158.58"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("D26").Value = "0.1350"
$ws.Range("E26").Value = "  -1.63%  "

$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").Value = "0.07264"
$ws.Range("E28").Value = "  +12.49%  "

$ws.Range("D29").Value = "1.462"
$ws.Range("E29").Value = "  +2.83%  "

$ws.Range("D30").Value = "1.481"
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").Value = "4.040"
$ws.Range("E31").Value = "  -1.33%  "

$ws.Range("E32").Value = "  -1.41%  "

$ws.Range("D33").Value = "1.824"
$ws.Range("E33").Value = "  -0.49%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").Value = "0.7091"
$ws.Range("E35").Value = "  +1.74%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").Value = "0.01841"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").Value = "6.926"
$ws.Range("E38").Value = "  +3.81%  "

$ws.Range("D39").Value = "2.817"
$ws.Range("E39").Value = "  -0.69%  "

$ws.Range("D40").Value = "1.234.78"
$ws.Range("E40").Value = "  -2.68%  "

$ws.Range("D41").Value = "0.9540"
$ws.Range("E41").Value = "  +4.38%  "

$ws.Range("E42").Value = "  +0.31%  "

$ws.Range("D43").Value = "2.003.19"
$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("D44").Value = "101.05"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").Value = "65.53"
$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "1.727"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000118"
$ws.Range("E47").Value = "  +2.18%  "

$ws.Range("D48").Value = "6.956"
$ws.Range("E48").Value = "  -1.79%  "

$ws.Range("D49").Value = "8.918"
$ws.Range("E49").Value = "  -1.32%  "

$ws.Range("D50").Value = "0.1138"
$ws.Range("E50").Value = "  -2.21%  "

$ws.Range("E51").Value = "  -1.58%  "
